# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" on all sheets
#    that reference it (Overview!E2:F2/E3:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2. Narrow the "Status" column(s) that held that text:
#    Overview columns E & F, zh-cn column C, de-de column C.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update the status text wherever it currently appears ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Narrow the status columns (Overview E:F, zh-cn C, de-de C) ---
$newColumnWidth = 13.4101845877511

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
